# cost_model.xlsx edit script
# - Rename CHAINSAW -> CHAINSAWold
# - Insert a new CHAINSAW sheet (cost roll-up) after CHAINSAWold
# - Rename defined names input1/input2 -> batteryCost/motorCost, add new names
# - Bar!B2 45 -> 30
# - Sheet1!B3/B4 new values
# - misc selection/active-sheet updates

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Rename old CHAINSAW sheet, insert new one after it
# ---------------------------------------------------------------------------
$oldChainsaw = $wb.Worksheets.Item("CHAINSAW")
$oldChainsaw.Name = "CHAINSAWold"
$oldChainsaw.Range("A10").Select() | Out-Null

$chainsaw = $wb.Worksheets.Add($null, $oldChainsaw)
$chainsaw.Name = "CHAINSAW"

Write-Host "done part 1"
